$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 274
$ws.Range("J2").Value = 1058
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 280
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 184
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 125
$ws.Range("T2").Value = 175
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 1625
$ws.Range("X2").Value = 1614
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 11
